$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "2023" column (Q) of data
$ws.Range("Q2").Value = $null
$ws.Range("Q3").Value = 2023
$ws.Range("Q4").Value = 279.01945525291825
$ws.Range("Q5").Value = 1792.7
$ws.Range("Q6").Value = 6425

# Copy formatting (styles) from column P onto the new column Q
$ws.Range("P2:P6").Copy()
$ws.Range("Q2:Q6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Adjust row heights for rows 4 and 5 as seen in the target workbook
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27.75

# Reset the selection to the default (A1) so no stray selection remains
$ws.Range("A1").Select()
